$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 180; $row -le 234; $row++) {
    $n = $row - 1
    $label = "B{0:D3}" -f $n
    $ws.Cells.Item($row, 1).Value = $label
    $ws.Cells.Item($row, 2).Value = 2018
    $ws.Cells.Item($row, 3).Value = 0
    $ws.Cells.Item($row, 4).Value = 0
    $ws.Cells.Item($row, 5).Value = 0
    $ws.Cells.Item($row, 6).Value = 0
    $ws.Cells.Item($row, 7).Value = 0
    $ws.Cells.Item($row, 8).Value = 0
}
